$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Remove the worker "LUIS HUMBERTO RUIZ DE LA CRUZ" (row 16) entirely - shifts
# the remaining table rows (JUVENAL x3 and MARIA PAULA) up by one.
$ws.Rows.Item(16).Delete()

# Re-order the "Periodo Mora" values for JUVENAL MARTINEZ CERVANTES so they
# read ascending (1910, 1911, 1912) instead of the previous descending order.
$ws.Range("E16").Value = "1910"
$ws.Range("E18").Value = "1912"

# Update the summary figures at the top of the account statement.
$ws.Range("E11").Value = 163200
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 4
